$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename the first sheet from "Tabelle1" to "Sheet1"
$ws1.Name = "Sheet1"

# Sheet1: B2 and B7 hold numbers that should become text ("1" / "2")
$ws1.Range("B2").NumberFormat = "@"
$ws1.Range("B2").Value = "1"

$ws1.Range("B7").NumberFormat = "@"
$ws1.Range("B7").Value = "2"

# meta sheet: B3 holds a number (2024) that should become text "2024"
$ws2.Range("B3").NumberFormat = "@"
$ws2.Range("B3").Value = "2024"

# Move the selection on the meta sheet from A4 to B3
$ws2.Range("B3").Select() | Out-Null

# Make Sheet1 the active tab (was "meta" before)
$ws1.Activate() | Out-Null
